$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 with same style (bold, border, centered) as existing header cells
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate time_taken values for rows 2-61
$values = @(
    "2021-10-05 13:42:05.081484",
    "2021-10-05 13:42:05.081498",
    "2021-10-05 13:42:05.081501",
    "2021-10-05 13:42:05.081505",
    "2021-10-05 13:42:05.081508",
    "2021-10-05 13:42:05.081511",
    "2021-10-05 13:42:05.081514",
    "2021-10-05 13:42:05.081517",
    "2021-10-05 13:42:05.081521",
    "2021-10-05 13:42:05.081524",
    "2021-10-05 13:42:05.081527",
    "2021-10-05 13:42:05.081530",
    "2021-10-05 13:42:05.081533",
    "2021-10-05 13:42:05.081536",
    "2021-10-05 13:42:05.081539",
    "2021-10-05 13:42:05.081542",
    "2021-10-05 13:42:05.081545",
    "2021-10-05 13:42:05.081548",
    "2021-10-05 13:42:05.081552",
    "2021-10-05 13:42:05.081555",
    "2021-10-05 13:42:05.081558",
    "2021-10-05 13:42:05.081561",
    "2021-10-05 13:42:05.081564",
    "2021-10-05 13:42:05.081567",
    "2021-10-05 13:42:05.081570",
    "2021-10-05 13:42:05.081574",
    "2021-10-05 13:42:05.081577",
    "2021-10-05 13:42:05.081580",
    "2021-10-05 13:42:05.081583",
    "2021-10-05 13:42:05.081586",
    "2021-10-05 13:42:05.081589",
    "2021-10-05 13:42:05.081592",
    "2021-10-05 13:42:05.081595",
    "2021-10-05 13:42:05.081599",
    "2021-10-05 13:42:05.081602",
    "2021-10-05 13:42:05.081605",
    "2021-10-05 13:42:05.081608",
    "2021-10-05 13:42:05.081611",
    "2021-10-05 13:42:05.081614",
    "2021-10-05 13:42:05.081617",
    "2021-10-05 13:42:05.081620",
    "2021-10-05 13:42:05.081623",
    "2021-10-05 13:42:05.081627",
    "2021-10-05 13:42:05.081630",
    "2021-10-05 13:42:05.081632",
    "2021-10-05 13:42:05.081635",
    "2021-10-05 13:42:05.081639",
    "2021-10-05 13:42:05.081642",
    "2021-10-05 13:42:05.081645",
    "2021-10-05 13:42:05.081648",
    "2021-10-05 13:42:05.081651",
    "2021-10-05 13:42:05.081654",
    "2021-10-05 13:42:05.081657",
    "2021-10-05 13:42:05.081660",
    "2021-10-05 13:42:05.081663",
    "2021-10-05 13:42:05.081666",
    "2021-10-05 13:42:05.081669",
    "2021-10-05 13:42:05.081672",
    "2021-10-05 13:42:05.081675",
    "2021-10-05 13:42:05.081678"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

